$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    3  = -11.1511
    14 = -13.65879999999999
    21 = -12.5659
    23 = -12.36690000000001
    25 = -13.6239
    26 = -12.62110000000001
    29 = -10.99670000000001
    53 = -10.41890000000001
    57 = -13.94379999999999
    59 = -12.69580000000001
    69 = -11.5998
    79 = -10.82140000000001
    83 = -13.5599
    91 = -10.2832
    93 = -11.1865
}

foreach ($row in $updates.Keys) {
    $ws.Range("C$row").Value = $updates[$row]
}
